$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: columns A/B/C/D/F/G hold free-form text in the source workbook
# (numbers/dates included) that must stay text rather than being
# auto-coerced to a number/date serial by Excel on assignment, so each
# text write is given a leading apostrophe (Excel's classic "force text"
# marker, stripped from the stored value). Column E is the one true
# numeric column and is written as a plain number.

$ws.Range("A25").Value = "'05-12-2025"
$ws.Range("B25").Value = "'010965012-Medha Sub Division Office Coll."
$ws.Range("D25").Value = "'2025-12-05"
$ws.Range("E25").Value = 28180
$ws.Range("F25").Value = "'"
$ws.Range("B26").Value = "'020965018-Kai Lalsingrao Shinde Gr.Big.Sheti Sah.Pat.Ltd. Br. Medha"
$ws.Range("E26").Value = 32740
$ws.Range("B27").Value = "'020965021-KAI.LALSINGRAO BAPUSO SHINDE SAH.PAT.LTD.,KUDAL, BR.KARAHAR"
$ws.Range("E27").Value = 4020
$ws.Range("B28").Value = "'020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"
$ws.Range("C28").Value = "'Cheque"
$ws.Range("D28").Value = "'"
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = "'रद्द झालेल्या दोन रिसिटची रक्कम रु 2020.00 एवढी बँकेमार्फत नजरचुकीने भरली गेल्याने, रु 2020.00  एवढी रक्कम दिनांक 06.12.2025 च्या एकूण कलेक्शन मधून कमी भरणार असल्याचे सांगण्यात आले. ( आज रोजीचे एकूण कलेक्शन 57110.00 असे आहे )"
$ws.Range("C29").Value = "'NEFT"
$ws.Range("C30").Value = "'Total"
$ws.Range("C31").Value = "'Cash"
$ws.Range("D31").Value = "'2025-12-05"
$ws.Range("E31").Value = 59730
$ws.Range("A33").Value = "'06-12-2025"
$ws.Range("C33").Value = "'Cheque"
$ws.Range("D33").Value = "'"
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = "'दिनांक 05.12.2025 रोजी रद्द झालेल्या दोन रिसिटची रक्कम रु 2020.00 एवढी बँकेमार्फत नजरचुकीने भरली गेल्याने, रु 2020.00 एवढी रक्कम दिनांक 06.12.2025 च्या एकूण कलेक्शन मधून कमी भरलेली आहे."
$ws.Range("C34").Value = "'NEFT"
$ws.Range("C35").Value = "'Total"
$ws.Range("C36").Value = "'Cash"
$ws.Range("D36").Value = "'2025-12-06"
$ws.Range("E36").Value = 47860
$ws.Range("G36").Value = "'2025-12-24"

# Remove the now-obsolete trailing Total row (old row 37); this also
# shrinks the sheet dimension from G37 to G36, matching the target.
$ws.Rows.Item(37).Delete()
